$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the formatting of the "YouTube Channel" hyperlink cell (G2) into a
# scratch cell far outside the table, so it survives the upcoming column
# deletions. We'll use it afterwards to restore the exact font/style, since
# Hyperlinks.Add() resets a cell to the generic built-in "Hyperlink" style.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("Z100").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Hyperlinks don't follow their cells when columns are inserted/deleted in
# this object model, so drop them all now and re-create the single survivor
# later at its new address.
$ws.Hyperlinks.Delete()

# Remove the columns that are dropped from the table:
#   H - "آخرین به‌روزرسانی اطلاعات"
#   F - "تاریخ آخرین قسمت منتشر شده"
#   E - "نماد" (the logo column)
# Column G ("لینک" / YouTube Channel) shifts left and becomes column E.
$ws.Columns("H").Delete()
$ws.Columns("F").Delete()
$ws.Columns("E").Delete()

# Rename the remaining header from "لینک" to "کانال منتشرکننده".
$ws.Range("E1").Value = "کانال منتشرکننده"

# Re-create the hyperlink on (the new) E2. Leave TextToDisplay unset so the
# existing cell text ("YouTube Channel", carried over from old G2) is kept
# as-is instead of being overwritten.
$ws.Hyperlinks.Add($ws.Range("E2"), "https://youtube.com/@EsmNadareh?si=IFvIYOsNRxbD9af0") | Out-Null

# Restore the original font/style (Shabnam 16, underlined hyperlink colour)
# that Hyperlinks.Add just reset, from the stash made earlier. Three columns
# were removed since the stash was made, so it now lives at W100 (Z - 3).
$ws.Range("W100").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("W100").Clear() | Out-Null

# Keep the "fit to page" print attributes alive across the edit.
$ws.PageSetup.FitToPagesWide = 0
$ws.PageSetup.FitToPagesTall = 0

# Update the selection/view: no more scrolled-to-C1 view, select B3 instead.
$ws.Range("B3").Select() | Out-Null
